# Refresh the "cryptos" price list: Price (D) and Volume(1h) (E) columns
# for each coin row. Numeric-looking price strings (e.g. "554.57") are
# written with a temporary Text number format so Excel keeps them as
# literal strings (matching the source data's inline-string cells)
# instead of silently coercing them to floating-point numbers; the
# format is cleared right after so the cell's style stays untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.194.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.576.22"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.11%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.76"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.24%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.582.50"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.72"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  +11.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.350"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.029.47"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.175.94"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.98"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.96%  "
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.578.70"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.32%  "
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.85"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.60"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.471"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -3.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.38"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0771"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.49%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.69"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.06"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.17"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.28"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.851"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.48"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "288.69"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.29"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.86%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0973"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0529"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.62"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.936.55"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.27%  "
